$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Register the hyperlink destination for B6 first (its style gets reset by Add,
# so the cell formatting/value is (re)applied afterwards to match the other rows).
$url = "http://raspberrypihq.com/how-to-share-a-folder-with-a-windows-computer-from-a-raspberry-pi/"
$ws.Hyperlinks.Add($ws.Range("B6"), $url)

# Copy formatting from the row above (A4/B4) so row 6 reuses the existing
# "plain text" and "Hyperlink" cell styles instead of creating new ones.
$ws.Range("A4").Copy($ws.Range("A6"))
$ws.Range("B4").Copy($ws.Range("B6"))

# Fill in the new PDO row: description text in A6, link text in B6.
# Set B6 before A6 so the shared-string table gets the URL at the lower index,
# matching the order new strings were appended upstream.
$ws.Range("B6").Value = $url
$ws.Range("A6").Value = "Setup samba share"

# New column B gets an explicit width, as seen in the target sheet.
$ws.Columns.Item(2).ColumnWidth = 38.8

# Sheet dimension now extends to row 6 and the active selection moves there too.
$ws.Range("A6").Select()
